{"js": "// Minor updates: \"Chapter 2\" -> \"Chapter 3\" in the title, and\n// \"wont\" -> \"won't\" (typographic apostrophe) in the second paragraph.\n// The editing Word session also leaves its \"last edit\" marker (the\n// hidden _GoBack bookmark) at the site of the final keystroke, i.e.\n// right after the newly-typed apostrophe, so we relocate it there too.\n\nconst body = context.document.body;\n\n// The hidden \"_GoBack\" bookmark tracks the location of the most recent\n// edit. Remove the old one first (bookmark names must stay unique) -\n// it gets re-created below at the new edit location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 1. Title: \"Chapter 2\" -> \"Chapter 3\"\nconst titleMatches = body.search(\"Chapter 2\", { matchCase: true, matchWholeWord: false });\ntitleMatches.load(\"items\");\nawait context.sync();\n\nif (titleMatches.items.length > 0) {\n  titleMatches.items[0].insertText(\"Chapter 3\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2. \"wont\" -> \"won't\" (curly/typographic apostrophe U+2019)\nconst wontMatches = body.search(\"wont\", { matchCase: true, matchWholeWord: true });\nwontMatches.load(\"items\");\nawait context.sync();\n\nif (wontMatches.items.length > 0) {\n  const wontRange = wontMatches.items[0];\n  wontRange.insertText(\"won\\u2019t\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Re-find a unique anchor right at the final \"t\" (the word now reads\n  // \"won\\u2019t need to worry...\") so the bookmark can be dropped exactly\n  // between the apostrophe and the closing \"t\", matching where the\n  // cursor sat after typing the apostrophe.\n  const anchorMatches = body.search(\"t need to worry\", { matchCase: true });\n  anchorMatches.load(\"items\");\n  await context.sync();\n\n  if (anchorMatches.items.length > 0) {\n    const rightBeforeFinalT = anchorMatches.items[0].getRange(\"Start\");\n    rightBeforeFinalT.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Minor updates: \"Chapter 2\" -> \"Chapter 3\" in the title, and\n# \"wont\" -> \"won't\" (typographic apostrophe) in the second paragraph.\n# The editing Word session also leaves its \"last edit\" marker (the\n# hidden _GoBack bookmark) at the site of the final keystroke, i.e.\n# right after the newly-typed apostrophe, so we relocate it there too.\n\n$d = $word.ActiveDocument\n\n# The hidden \"_GoBack\" bookmark tracks the location of the most recent\n# edit. Remove the old one first (bookmark names must stay unique) -\n# it gets re-created below at the new edit location.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 1. Title: \"Chapter 2\" -> \"Chapter 3\"\n$titleFind = $d.Content.Find\n$titleFind.ClearFormatting()\n$titleFind.Text = \"Chapter 2\"\n$titleFind.MatchCase = $true\n$titleFind.MatchWholeWord = $false\n$titleFind.Execute()\nif ($titleFind.Found) {\n    $titleFind.Parent.Text = \"Chapter 3\"\n}\n\n# 2. \"wont\" -> \"won't\" (curly/typographic apostrophe U+2019)\n$wontFind = $d.Content.Find\n$wontFind.ClearFormatting()\n$wontFind.Text = \"wont\"\n$wontFind.MatchCase = $true\n$wontFind.MatchWholeWord = $true\n$wontFind.Execute()\nif ($wontFind.Found) {\n    $wontRange = $wontFind.Parent\n    $wontRange.Text = \"won\" + [char]0x2019 + \"t\"\n}\n\n# Drop the \"_GoBack\" bookmark exactly between the apostrophe and the\n# closing \"t\", matching where the cursor sat after typing the apostrophe.\n$anchorFind = $d.Content.Find\n$anchorFind.ClearFormatting()\n$anchorFind.Text = \"t need to worry\"\n$anchorFind.MatchCase = $true\n$anchorFind.Execute()\nif ($anchorFind.Found) {\n    $anchorRange = $anchorFind.Parent\n    $anchorRange.Collapse(1)\n    $d.Bookmarks.Add(\"_GoBack\", $anchorRange)\n}\n"}
